$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3161.9092
$ws.Range("J19").Value = 3141
$ws.Range("L19").Value = 3141
$ws.Range("N19").Value = -3491

$ws.Range("H62").Value = 4381.75
$ws.Range("I62").Value = 4381.75
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4381.75
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3757.75
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4381.75
$ws.Range("I65").Value = 4381.75
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 21908.75
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -18788.75
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 5378.2
$ws.Range("I74").Value = 5042.5557
$ws.Range("K74").Value = 5042.5557
$ws.Range("M74").Value = -4106.5557

$ws.Range("H77").Value = 5378.2
$ws.Range("I77").Value = 5042.5557
$ws.Range("K77").Value = 25212.7785
$ws.Range("M77").Value = -20532.7785

$ws.Range("H98").Value = 1974.75
$ws.Range("J98").Value = 1392.5
$ws.Range("L98").Value = 1392.5
$ws.Range("N98").Value = -4388.5

$ws.Range("H122").Value = 1974.75
$ws.Range("J122").Value = 1392.5
$ws.Range("L122").Value = 4177.5
$ws.Range("N122").Value = -9077.5

$ws.Range("H132").Value = 66099.625
$ws.Range("I132").Value = 80463.69500000001
$ws.Range("K132").Value = 241391.085
$ws.Range("M132").Value = -238861.085

$ws.Range("H137").Value = 1849.44
$ws.Range("I137").Value = 1880.4
$ws.Range("K137").Value = 5641.200000000001
$ws.Range("M137").Value = -3091.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4388.4375
$ws.Range("I32").Value = 4506
$ws.Range("J32").Value = 2625
$ws.Range("K32").Value = 4506
$ws.Range("L32").Value = 2625
$ws.Range("M32").Value = -4219
$ws.Range("N32").Value = -3199

$ws.Range("H74").Value = 2932.1052
$ws.Range("J74").Value = 1942
$ws.Range("L74").Value = 1942
$ws.Range("N74").Value = -3690

$ws.Range("H77").Value = 2932.1052
$ws.Range("J77").Value = 1942
$ws.Range("L77").Value = 9710
$ws.Range("N77").Value = -18446

$ws.Range("H108").Value = 44999
$ws.Range("J108").Value = 44999
$ws.Range("L108").Value = 44999
$ws.Range("N108").Value = -52679

$ws.Range("H132").Value = 18522684
$ws.Range("I132").Value = 3204.0454
$ws.Range("K132").Value = 9612.136200000001
$ws.Range("M132").Value = -7082.136200000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 7964.9
$ws.Range("I3").Value = 11099.182
$ws.Range("J3").Value = 4134.1113
$ws.Range("K3").Value = 11099.182
$ws.Range("L3").Value = 4134.1113
$ws.Range("M3").Value = -10985.182
$ws.Range("N3").Value = -4362.1113

$ws.Range("H40").Value = 44929
$ws.Range("J40").Value = 44929
$ws.Range("L40").Value = 44929
$ws.Range("N40").Value = -45459

$ws.Range("H96").Value = 35455.2
$ws.Range("I96").Value = 35455.2
$ws.Range("K96").Value = 35455.2
$ws.Range("M96").Value = -32709.2

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("J134").Value = 2000
$ws.Range("L134").Value = 6000
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2537.3044
$ws.Range("I31").Value = 2657.5
$ws.Range("J31").Value = 2406.182
$ws.Range("K31").Value = 2657.5
$ws.Range("L31").Value = 2406.182
$ws.Range("M31").Value = -2362.5
$ws.Range("N31").Value = -2996.182

$ws.Range("H34").Value = 2537.3044
$ws.Range("I34").Value = 2657.5
$ws.Range("J34").Value = 2406.182
$ws.Range("K34").Value = 2657.5
$ws.Range("L34").Value = 2406.182
$ws.Range("M34").Value = -2455.5
$ws.Range("N34").Value = -2810.182

$ws.Range("H99").Value = 3305.1428
$ws.Range("J99").Value = 3875.3333
$ws.Range("L99").Value = 3875.3333
$ws.Range("N99").Value = -6871.3333

$ws.Range("H109").Value = 43998
$ws.Range("J109").Value = 43998
$ws.Range("L109").Value = 43998
$ws.Range("N109").Value = -46078

$ws.Range("H126").Value = 3305.1428
$ws.Range("J126").Value = 3875.3333
$ws.Range("L126").Value = 11625.9999
$ws.Range("N126").Value = -16565.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 19857.143
$ws.Range("J105").Value = 19857.143
$ws.Range("L105").Value = 59571.429
$ws.Range("N105").Value = -64813.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 21332.666
$ws.Range("J98").Value = 21332.666
$ws.Range("L98").Value = 21332.666
$ws.Range("N98").Value = -27322.666

$ws.Range("H102").Value = 2488.2
$ws.Range("I102").Value = 2488.2
$ws.Range("K102").Value = 2488.2
$ws.Range("M102").Value = -866.1999999999998

$ws.Range("H122").Value = 5292.7646
$ws.Range("I122").Value = 6208.8335
$ws.Range("J122").Value = 3094.2
$ws.Range("K122").Value = 18626.5005
$ws.Range("L122").Value = 9282.599999999999
$ws.Range("M122").Value = -16176.5005
$ws.Range("N122").Value = -14182.6

$ws.Range("H132").Value = 4497.0713
$ws.Range("I132").Value = 4363.6665
$ws.Range("J132").Value = 4737.2
$ws.Range("K132").Value = 13090.9995
$ws.Range("L132").Value = 14211.6
$ws.Range("M132").Value = -10560.9995
$ws.Range("N132").Value = -19271.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7411.4287
$ws.Range("J22").Value = 7411.4287
$ws.Range("L22").Value = 7411.4287
$ws.Range("N22").Value = -8001.4287

$ws.Range("H27").Value = 7411.4287
$ws.Range("J27").Value = 7411.4287
$ws.Range("L27").Value = 7411.4287
$ws.Range("N27").Value = -7625.4287

$ws.Range("H46").Value = 4960
$ws.Range("I46").Value = 4399.7144
$ws.Range("K46").Value = 4399.7144
$ws.Range("M46").Value = -4211.7144

$ws.Range("H82").Value = 2565.7693
$ws.Range("J82").Value = 3663.75
$ws.Range("L82").Value = 3663.75
$ws.Range("N82").Value = -4385.75

$ws.Range("H85").Value = 2565.7693
$ws.Range("J85").Value = 3663.75
$ws.Range("L85").Value = 3663.75
$ws.Range("N85").Value = -6159.75

$ws.Range("H122").Value = 7738.1177
$ws.Range("J122").Value = 8591.666999999999
$ws.Range("L122").Value = 25775.001
$ws.Range("N122").Value = -30675.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 18887.125
$ws.Range("I81").Value = 21182.834
$ws.Range("K81").Value = 42365.668
$ws.Range("M81").Value = -41304.668

$ws.Range("H84").Value = 18887.125
$ws.Range("I84").Value = 21182.834
$ws.Range("K84").Value = 211828.34
$ws.Range("M84").Value = -206524.34

$ws.Range("H115").Value = 49947
$ws.Range("J115").Value = 49947
$ws.Range("L115").Value = 49947
$ws.Range("N115").Value = -53081

$ws.Range("H132").Value = 3182.3333
$ws.Range("I132").Value = 3217.5715
$ws.Range("K132").Value = 9652.7145
$ws.Range("M132").Value = -7122.7145
